$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing score values (set to 5) for the specified cells
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 5

$ws.Range("E16").Value = 5
$ws.Range("F16").Value = 5

$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 5
$ws.Range("F17").Value = 5

$ws.Range("E20").Value = 5
$ws.Range("F20").Value = 5

$ws.Range("E26").Value = 5
$ws.Range("F26").Value = 5

$ws.Range("G31").Value = 5

# Update the active cell selection to match the saved view state
$ws.Range("F6").Select()
